$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that used to sit right
#    after the H1 title at the top of the document.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2. Insert a new bold paragraph, "Play Day and Night Slot Free | Stunning
#    Graphics & Unique Features", right before the final (italic) paragraph
#    near the bottom of the document.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPoint = $lastPara.Range.Start
$insertRange = $d.Range($insertPoint, $insertPoint)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Day and Night Slot Free | Stunning Graphics &amp; Unique Features</w:t></w:r></w:p>'
$insertRange.InsertXML($newParaXml)

# The XML fragment above merges into the start of the existing last
# paragraph, so split it into its own paragraph right after the inserted
# text.
$newTextLength = 66
$splitPoint = $insertPoint + $newTextLength
$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 3. Replace the body text of the final (italic) paragraph with the new
#    meta-description copy, keeping its existing italic run formatting.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$finalPara = $d.Paragraphs.Item($count)
$oldText = 'Create a feature image fitting the game "Day and Night": Please create an image in cartoon style featuring a happy Maya warrior with glasses, surrounded by the opposing forces of light and darkness. The warrior should be holding a staff or other ancient Egyptian weapon, and standing confidently amidst the clash of the two deities, Ra and Bastet. The image should be colorful and dynamic, with radiant orange hues on one side and a dark, mystical blue on the other. The name of the game, "Day and Night" should be featured prominently in the image, and any additional text or graphics should be inspired by ancient Egyptian mythology. Overall, the image should capture the excitement and intrigue of this unique and visually stunning online slot game.'
$newText = 'Read our Day and Night slot review and play for free. Stunning graphics, free spins and low variance make it a game for all players.'
$finalPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
